# Update cryptocurrency price/volume data in the worksheet to reflect the
# latest scrape values (commit: 'Updated cryptos list ... with GitHub Actions').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text would otherwise be auto-parsed by Excel
# as a number (e.g. '276.86' -> 276.86, '1.000' -> 1, '0.06667' -> 6.667E-02),
# losing the original formatting (trailing zeros, exact digit grouping, etc.).
# Force these cells to Text format first so the literal string is preserved,
# exactly as it already is for the two-dot values like '25.785.47'.
$textFormatCells = @(
    "D5",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D15",
    "D16",
    "D17",
    "D19",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New cell values: column D (Price) and column E (Volume(1h)) updates.
$cellValues = @{
    "D2" = '25.785.47'
    "E2" = '  -3.93%  '
    "D3" = '1.815.52'
    "E3" = '  -3.13%  '
    "E4" = '  -0.09%  '
    "D5" = '276.86'
    "E5" = '  -8.26%  '
    "E6" = '  -0.14%  '
    "D7" = '0.5107'
    "E7" = '  -5.19%  '
    "E8" = '  -6.54%  '
    "D9" = '44.65'
    "E9" = '  -2.18%  '
    "D10" = '0.06667'
    "E10" = '  -7.26%  '
    "D11" = '20.01'
    "E11" = '  -7.25%  '
    "D12" = '0.8304'
    "D13" = '0.07866'
    "E13" = '  -3.57%  '
    "D14" = '1.834.15'
    "E14" = '  -2.17%  '
    "D15" = '5.078'
    "E15" = '  -3.54%  '
    "D16" = '87.42'
    "E16" = '  -6.37%  '
    "D17" = '1.000'
    "E17" = '  -0.07%  '
    "E18" = '  -4.23%  '
    "D19" = '0.000008027'
    "E19" = '  -6.13%  '
    "E20" = '  -0.08%  '
    "D21" = '25.858.35'
    "E21" = '  -3.85%  '
    "D22" = '4.725'
    "D23" = '10.01'
    "E23" = '  -6.39%  '
    "D24" = '6.077'
    "E24" = '  -4.90%  '
    "D25" = '141.17'
    "E25" = '  -4.15%  '
    "D26" = '2.190'
    "E26" = '  -3.10%  '
    "D27" = '1.673'
    "E27" = '  -3.79%  '
    "D28" = '17.07'
    "E28" = '  -5.35%  '
    "D29" = '109.48'
    "E29" = '  -4.12%  '
    "E30" = '  -7.96%  '
    "D31" = '4.234'
    "E31" = '  -7.95%  '
    "D32" = '0.08803'
    "E32" = '  -3.92%  '
    "D33" = '0.04876'
    "E33" = '  -2.04%  '
    "D34" = '0.7322'
    "E34" = '  -9.25%  '
    "D35" = '1.135'
    "E35" = '  -3.07%  '
    "E36" = '  -3.71%  '
    "D37" = '3.151'
    "E37" = '  -1.19%  '
    "E38" = '  -0.09%  '
    "D39" = '2.370'
    "E39" = '  -8.54%  '
    "D40" = '0.5213'
    "E40" = '  -13.77%  '
    "D41" = '0.01848'
    "E41" = '  -5.56%  '
    "D42" = '0.9562'
    "E42" = '  -11.00%  '
    "D43" = '111.71'
    "E43" = '  -3.56%  '
    "D44" = '6.197'
    "E44" = '  -5.93%  '
    "D45" = '8.030'
    "E45" = '  -9.21%  '
    "E46" = '  -0.11%  '
    "D47" = '0.4560'
    "E47" = '  -11.60%  '
    "E48" = '  -8.89%  '
    "D49" = '36.72'
    "D50" = '9.223'
    "D51" = '1.501'
    "E51" = '  -8.21%  '
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
